$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B7").Value  = "Outdoor Horticulture/Arable Farming"
$ws.Range("B10").Value = "Other (Agriculture)"
$ws.Range("C10").Value = "This balances energy in according to New Zealands Energy Balances from MBIE."
$ws.Range("B15").Value = "Other (Commercial)"
$ws.Range("C23").Value = "Wind grows faster in T$([char]0x16B)$([char]0x12B) as the Aluminium smelter exit in Kea results in less electricity demand for Kea. "
$ws.Range("C24").Value = "Large decarbonisation occurs in the industrial sector with most of the remaining fossil fuels in hard to abate sectors. Note - Only energy related emissions are included in the TIMES-NZ model. Emissions from Feedstock are not expressed"
$ws.Range("C25").Value = "In Kea, aluminium smelting reaches zero in 2030, and in T$([char]0x16B)$([char]0x12B) the demand stays as the smelter may continue to run or is replaced by another industry."
$ws.Range("B28").Value = "Fabricated Metal Product, Transport Equipment, Machinery and Equipment Manufacturing"
$ws.Range("C28").Value = "T$([char]0x16B)$([char]0x12B) assumes a higher growth rate than in Kea, because Kea tries to move away from emission-heavy manufacturing."
$ws.Range("C31").Value = "Kea assumes methanol production exits at 2032.  In T$([char]0x16B)$([char]0x12B) it exits in 2047. "
$ws.Range("B34").Value = "Other (Industry)"
$ws.Range("C44").Value = "A strong drive for electrification reduces emissions significantly. T$([char]0x16B)$([char]0x12B) sees greater use of hybrids in the short term, while Kea pushes hard on EVs."
$ws.Range("C47").Value = " In each time period, the share of EVs is larger in Kea than in T$([char]0x16B)$([char]0x12B) because the Kea scenario assumes there is a larger ability to access EVs.  Note - the end use demand for Road Transport is measured as Distance Travelled in the metric dropdown tool."

$ws.Range("C10").Select()
